# Jagadeesha Suchith.xlsx — "complate!!-> scrapping whole ipl"
#
# 1. Rename the sheet from the default "Sheet1" to the player's name.
# 2. Insert a new first column ("matchNo") in front of the existing data,
#    shifting teamName..result from A:L to B:M.
# 3. Populate the new column's header and the single data row's value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A — existing columns A:L shift right to B:M.
$ws.Columns.Item(1).Insert()

# New leading column: matchNo
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "20th"

# Rename the worksheet tab to match the player.
$ws.Name = "Jagadeesha Suchith"
